# edit.ps1 - apply the diff described:
#  1. Update the datetimeFigureOut field text from 21.10.2024 to 22.10.2024
#     everywhere it occurs (slide master, all slide layouts, notes master).
#  2. Add a new "Down Arrow" autoshape to slide 1.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update the date placeholders (datetimeFigureOut fields) wherever
#    they currently read "21.10.2024" -> "22.10.2024".
# ---------------------------------------------------------------------

function Update-DateFields {
    param($shapes)

    foreach ($shp in $shapes) {
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq "21.10.2024") {
                    $tr.Text = "22.10.2024"
                }
            }
        }
    }
}

# Slide Master
Update-DateFields -shapes $p.SlideMaster.Shapes

# All Slide Layouts
foreach ($layout in $p.SlideMaster.CustomLayouts) {
    Update-DateFields -shapes $layout.Shapes
}

# Notes Master
Update-DateFields -shapes $p.NotesMaster.Shapes

# ---------------------------------------------------------------------
# 2. Add the new down-arrow shape to slide 1.
# ---------------------------------------------------------------------

$s = $p.Slides.Item(1)

$msoShapeDownArrow = 56

$shape = $s.Shapes.AddShape($msoShapeDownArrow, 459.13, 150.54, 17.35, 146.57)
$shape.Name = "Pfeil: nach unten 227"
$shape.Rotation = 46.791383333333336
